$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rewrite the opening paragraph (FirstParagraph style): drop the
#    "Starting at age 13..." sentence, drop "up to", and drop the quoted
#    'trial by fire' aside.
# ---------------------------------------------------------------------------
$introFind = $d.Content
$introFind.Find.Execute("Starting at age 13", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$introFind.Expand(4) | Out-Null
$introRange = $d.Range($introFind.Start, $introFind.End - 1)
$introRange.Text = "I have managed remote teams of 12 engineers, designers and product managers in diverse industry sectors. While designing innovative products for clients, I" + [char]8217 + "ve learned how to be a lean and proficient product manager, facilitating efficient coordination among all team members."

# ---------------------------------------------------------------------------
# 2. "I'm competent in leading..." -> "I'm awesome at leading..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("I" + [char]8217 + "m competent in leading product teams", $true, $false, $false, $false, $false, $true, 1, $false, "I" + [char]8217 + "m awesome at leading product teams", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Remove the "I'm currently based in Ottawa..." paragraph entirely.
# ---------------------------------------------------------------------------
$ottawaFind = $d.Content
$ottawaFind.Find.Execute("I" + [char]8217 + "m currently based in Ottawa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ottawaFind.Expand(4) | Out-Null
$ottawaFind.Delete()

# ---------------------------------------------------------------------------
# 4. Move the "Education" section (Heading2 + DefinitionTerm + Definition)
#    from before "Professional Experience" to just after the Professional
#    Experience entries (right before "Other Contributions and Projects").
# ---------------------------------------------------------------------------
$destFind = $d.Content
$destFind.Find.Execute("Other Contributions and Projects", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$destPara = $destFind.Paragraphs.Item(1)
$insertPoint = $destPara.Range.Start

$eduHeadingText = "Education"
$eduTermText = "2002-2007"
$eduDefText = "B.Eng Systems Engineering Carleton University (Ottawa)"

$insertRange = $d.Range($insertPoint, $insertPoint)
$insertRange.InsertBefore($eduHeadingText + [char]13 + $eduTermText + [char]13 + $eduDefText + [char]13)

$p1Start = $insertPoint
$p1End = $p1Start + $eduHeadingText.Length
$p2Start = $p1End + 1
$p2End = $p2Start + $eduTermText.Length
$p3Start = $p2End + 1
$p3End = $p3Start + $eduDefText.Length

$eduHeadingRange = $d.Range($p1Start, $p1End)
$eduTermRange = $d.Range($p2Start, $p2End)
$eduDefRange = $d.Range($p3Start, $p3End)

$eduTermRange.Paragraphs.Item(1).Style = "DefinitionTerm"
$eduDefRange.Paragraphs.Item(1).Style = "Definition"

$d.Bookmarks.Add("education", $eduHeadingRange) | Out-Null

$eduBoldRange = $d.Range($p3Start, $p3Start + "B.Eng Systems Engineering".Length)
$eduBoldRange.Bold = 1

# Now delete the original Education section near the top of the resume.
$origFind = $d.Content
$origFind.Find.Execute("Education", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$origHeadingPara = $origFind.Paragraphs.Item(1)
$paras = $d.Paragraphs
$origIndex = 0
$k = 1
foreach ($pp in $paras) {
    if ($pp.Range.Start -eq $origHeadingPara.Range.Start) {
        $origIndex = $k
    }
    $k = $k + 1
}
$paras2 = $d.Paragraphs
$origP1 = $paras2.Item($origIndex)
$origP3 = $paras2.Item($origIndex + 2)
$origDelRange = $d.Range($origP1.Range.Start, $origP3.Range.End)
$origDelRange.Delete()

# ---------------------------------------------------------------------------
# 5. "Lead the product redesign" -> "Lead the new product design"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Lead the product redesign of our internal data exploration", $true, $false, $false, $false, $false, $true, 1, $false, "Lead the new product design of our internal data exploration", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. "core company KPIs in Metabase" -> "core company KPI reporting in Metabase"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Set up and implemented the core company KPIs in Metabase", $true, $false, $false, $false, $false, $true, 1, $false, "Set up and implemented the core company KPI reporting in Metabase", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7. "React" -> "React and Javascript" (standalone bullet item)
# ---------------------------------------------------------------------------
$reactRange = $d.Content
$reactRange.Find.Execute("React", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$reactRange.Text = "React and Javascript"

# ---------------------------------------------------------------------------
# 8. Willowbee.ca title change.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Willowbee.ca - Product Manager, Technical Manager, and Marketing. Founder", $true, $false, $false, $false, $false, $true, 1, $false, "Willowbee.ca - Engineering Manager and Product Manager. Co-Founder", 2) | Out-Null

# ---------------------------------------------------------------------------
# 9. "Will and Testament builder" -> "Will and Testament application"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("most advanced Canadian Will and Testament builder", $true, $false, $false, $false, $false, $true, 1, $false, "most advanced Canadian Will and Testament application", 2) | Out-Null

Write-Host "All edits applied"
